# This workbook logs daily Repollo (cabbage) price records for
# "Femacal de La Calera". Two new daily records were inserted at the top
# of the data block (rows 352-353), pushing all the existing records
# down by two rows (old row N -> new row N+2). The sheet's used range
# grows from A1:R442 to A1:R444 as a result.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the first data row that needs to
# shift (row 352), which pushes rows 352:442 down to 354:444 and grows
# the sheet dimension accordingly.
$ws.Rows("352:353").Insert()

# Row 352 - new record (2021-12-21, Primera)
$ws.Range("A352").Value = 3
$ws.Range("B352").Value = 'Femacal de La Calera'
$ws.Range("C352").Value = 'Coquimbo'
$ws.Range("D352").Value = 44551
$ws.Range("E352").Value = 5
$ws.Range("F352").Value = 100112006
$ws.Range("G352").Value = 'Repollo'
$ws.Range("H352").Value = 'Crespo record'
$ws.Range("I352").Value = 'Primera'
$ws.Range("J352").Value = 3200
$ws.Range("K352").Value = 600
$ws.Range("L352").Value = 700
$ws.Range("M352").Value = 650
$ws.Range("N352").Value = '$/unidad'
$ws.Range("O352").Value = 'Provincia de Quillota'
$ws.Range("P352").Value = 650
$ws.Range("Q352").Value = 1
$ws.Range("R352").Value = 'Hortaliza'

# Row 353 - new record (2021-12-21, Segunda)
$ws.Range("A353").Value = 3
$ws.Range("B353").Value = 'Femacal de La Calera'
$ws.Range("C353").Value = 'Coquimbo'
$ws.Range("D353").Value = 44551
$ws.Range("E353").Value = 5
$ws.Range("F353").Value = 100112006
$ws.Range("G353").Value = 'Repollo'
$ws.Range("H353").Value = 'Crespo record'
$ws.Range("I353").Value = 'Segunda'
$ws.Range("J353").Value = 1500
$ws.Range("K353").Value = 500
$ws.Range("L353").Value = 500
$ws.Range("M353").Value = 500
$ws.Range("N353").Value = '$/unidad'
$ws.Range("O353").Value = 'Provincia de Quillota'
$ws.Range("P353").Value = 500
$ws.Range("Q353").Value = 1
$ws.Range("R353").Value = 'Hortaliza'
